# Generate Report for Handoff
# Adds two new tracked files (6b715619-...md and b5666c1b-...md) as rows 4 & 5
# on the "Overview", "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$uuid1 = "6b715619-f78b-47f4-b3a1-36725e4b3e47"
$uuid2 = "b5666c1b-bb9a-4410-a41a-5fa3b5e7f5fb"

$xlf1Zh = "6b715619-f78b-47f4-b3a1-36725e4b3e47.19dc7c7e9032a6539d84ba710032ed2bba911746.zh-cn.xlf"
$xlf2Zh = "b5666c1b-bb9a-4410-a41a-5fa3b5e7f5fb.e3faf1db2c8f5fd4afe6079e62b2200aff415251.zh-cn.xlf"
$xlf1De = "6b715619-f78b-47f4-b3a1-36725e4b3e47.19dc7c7e9032a6539d84ba710032ed2bba911746.de-de.xlf"
$xlf2De = "b5666c1b-bb9a-4410-a41a-5fa3b5e7f5fb.e3faf1db2c8f5fd4afe6079e62b2200aff415251.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview" - two new rows (4 & 5)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "$uuid1.md"
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Range("D4").Value = "2016-03-21 16:39:32"

$wsOverview.Range("A5").Value = "$uuid2.md"
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-03-21 16:39:32"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/19dc7c7e9032a6539d84ba710032ed2bba911746/e2e/$uuid1.md", "", "", "$uuid1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/e3faf1db2c8f5fd4afe6079e62b2200aff415251/e2e/$uuid2.md", "", "", "$uuid2.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn" - two new rows (4 & 5)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = "$uuid1.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = $xlf1Zh
$wsZhCn.Range("E4").Value = "2016-03-21 16:39:29"
$wsZhCn.Range("H4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("J4").Value = "Include"

$wsZhCn.Range("A5").Value = "$uuid2.md"
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = $xlf2Zh
$wsZhCn.Range("E5").Value = "2016-03-21 16:39:29"
$wsZhCn.Range("H5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("J5").Value = "Include"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/19dc7c7e9032a6539d84ba710032ed2bba911746/e2e/$uuid1.md", "", "", "$uuid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/19dc7c7e9032a6539d84ba710032ed2bba911746/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlf1Zh", "", "", $xlf1Zh)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/e3faf1db2c8f5fd4afe6079e62b2200aff415251/e2e/$uuid2.md", "", "", "$uuid2.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e3faf1db2c8f5fd4afe6079e62b2200aff415251/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlf2Zh", "", "", $xlf2Zh)

# ---------------------------------------------------------------------------
# Sheet "de-de" - two new rows (4 & 5)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = "$uuid1.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = $xlf1De
$wsDeDe.Range("E4").Value = "2016-03-21 16:39:32"
$wsDeDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("J4").Value = "Include"

$wsDeDe.Range("A5").Value = "$uuid2.md"
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = $xlf2De
$wsDeDe.Range("E5").Value = "2016-03-21 16:39:32"
$wsDeDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("J5").Value = "Include"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/19dc7c7e9032a6539d84ba710032ed2bba911746/e2e/$uuid1.md", "", "", "$uuid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/19dc7c7e9032a6539d84ba710032ed2bba911746/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlf1De", "", "", $xlf1De)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/e3faf1db2c8f5fd4afe6079e62b2200aff415251/e2e/$uuid2.md", "", "", "$uuid2.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e3faf1db2c8f5fd4afe6079e62b2200aff415251/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlf2De", "", "", $xlf2De)
